$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.015.19"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.622.93"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.01%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.48"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -1.40%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.94"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -0.65%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0839"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").Value = "1.848.12"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.624.09"
$ws.Range("E13").Value = "  -0.98%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.11"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "26.985.99"
$ws.Range("E16").Value = "  -0.66%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.25"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("E18").Value = "  -0.32%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.40"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("E20").Value = "  +0.07%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -2.05%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.34"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -7.42%  "
$ws.Range("E24").Value = "  -1.96%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.63"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.53%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.48"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  -3.88%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -1.23%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.29"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -2.70%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.702"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +27.71%  "
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "1.337.69"
$ws.Range("E35").Value = "  +2.76%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  -0.87%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.840"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  +0.01%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.797"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -1.01%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("E43").Value = "  +0.11%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.77"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").Value = "1.759.92"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("E47").Value = "  +1.88%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.835"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +13.75%  "
$ws.Range("E49").Value = "  +0.26%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0990"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +3.52%  "
$ws.Range("E51").Value = "  -1.43%  "
